$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.340.90'
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.724.90'
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.06'
$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("E6").Value = '  +3.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '657.61'
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("E8").Value = '  +2.88%  '

$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.07'
$ws.Range("E10").Value = '  -1.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.722.98'
$ws.Range("E11").Value = '  +1.49%  '

$ws.Range("E12").Value = '  +18.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.89'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.92'
$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.418.94'
$ws.Range("E16").Value = '  +1.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '97.101.15'
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.01'
$ws.Range("E18").Value = '  +0.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.721.99'
$ws.Range("E19").Value = '  +1.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.08'
$ws.Range("E20").Value = '  +2.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.80'
$ws.Range("E21").Value = '  -0.48%  '

$ws.Range("E22").Value = '  -4.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '525.37'
$ws.Range("E23").Value = '  -1.40%  '

$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("E25").Value = '  +11.25%  '

$ws.Range("E26").Value = '  -4.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '106.31'
$ws.Range("E27").Value = '  +3.88%  '

$ws.Range("E28").Value = '  +15.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.920.77'
$ws.Range("E29").Value = '  +1.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.53'
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("E32").Value = '  -0.90%  '

$ws.Range("E33").Value = '  -0.10%  '

$ws.Range("E34").Value = '  +3.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.85'
$ws.Range("E35").Value = '  -2.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.66'
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '643.39'
$ws.Range("E38").Value = '  -1.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.595'
$ws.Range("E39").Value = '  -0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.76'
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.167'
$ws.Range("E42").Value = '  +2.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.495'
$ws.Range("E43").Value = '  +11.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.76'
$ws.Range("E44").Value = '  -0.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.73'
$ws.Range("E45").Value = '  +4.84%  '

$ws.Range("E46").Value = '  +1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.970'
$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("E48").Value = '  -0.30%  '

$ws.Range("E49").Value = '  +2.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.63'
$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.68'
$ws.Range("E51").Value = '  -0.70%  '
